# Apply the numeric updates described by the diff to both the "展览"
# sheet and the "全部类型" sheet (which duplicates the same rows).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6284
$ws1.Range("F5").Value = 25
$ws1.Range("F8").Value = 1401
$ws1.Range("F9").Value = 93

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6284
$ws4.Range("F5").Value = 25
$ws4.Range("F12").Value = 1401
$ws4.Range("F13").Value = 93
